# Auto-generated: updates cryptos price/volume cells per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '44.406.55'
$ws.Range('E2').Value = '  +1.27%  '
$ws.Range('D3').Value = '2.221.80'
$ws.Range('E3').Value = '  -0.68%  '
$ws.Range('E4').Value = '  +0.29%  '
$ws.Range('D5').Value = "'303.05"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.30%  '
$ws.Range('D6').Value = "'89.88"
$ws.Range('D6').Style = 'Normal'
$ws.Range('D7').Value = "'0.555"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.66%  '
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('D9').Value = "'0.497"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -4.71%  '
$ws.Range('D10').Value = "'33.57"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.56%  '
$ws.Range('D11').Value = "'0.0782"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.06%  '
$ws.Range('D12').Value = "'6.95"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E13').Value = '  -0.54%  '
$ws.Range('D14').Value = '2.562.67'
$ws.Range('E14').Value = '  -0.57%  '
$ws.Range('D15').Value = '2.328.72'
$ws.Range('E15').Value = '  +0.26%  '
$ws.Range('D16').Value = "'0.804"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.43%  '
$ws.Range('D17').Value = "'13.16"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.83%  '
$ws.Range('D18').Value = '44.159.20'
$ws.Range('E18').Value = '  +0.46%  '
$ws.Range('D19').Value = '0.0₃0910'
$ws.Range('E19').Value = '  -5.29%  '
$ws.Range('D20').Value = "'6.02"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -4.90%  '
$ws.Range('D21').Value = "'11.33"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -5.39%  '
$ws.Range('D22').Value = "'64.31"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.70%  '
$ws.Range('D23').Value = "'233.45"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.12%  '
$ws.Range('E24').Value = '  -1.64%  '
$ws.Range('E25').Value = '  +0.00%  '
$ws.Range('E26').Value = '  -4.10%  '
$ws.Range('D27').Value = "'2.26"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.53%  '
$ws.Range('D28').Value = "'9.46"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.64%  '
$ws.Range('D29').Value = "'36.31"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -8.94%  '
$ws.Range('D30').Value = "'19.53"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.48%  '
$ws.Range('D31').Value = "'5.59"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.14%  '
$ws.Range('D32').Value = "'146.59"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.85%  '
$ws.Range('E33').Value = '  +1.03%  '
$ws.Range('D34').Value = "'0.0756"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.37%  '
$ws.Range('D35').Value = "'2.99"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.00%  '
$ws.Range('D36').Value = "'0.106"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.97%  '
$ws.Range('E37').Value = '  -3.22%  '
$ws.Range('D38').Value = "'1.78"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.13%  '
$ws.Range('D39').Value = "'14.65"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +4.25%  '
$ws.Range('D40').Value = "'3.22"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -6.54%  '
$ws.Range('D41').Value = "'3.64"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.46%  '
$ws.Range('D42').Value = "'0.0288"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.97%  '
$ws.Range('E43').Value = '  +0.00%  '
$ws.Range('D44').Value = '1.769.36'
$ws.Range('E44').Value = '  +3.26%  '
$ws.Range('D45').Value = "'1.71"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +6.88%  '
$ws.Range('D46').Value = "'78.94"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.93%  '
$ws.Range('D47').Value = "'0.181"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -5.11%  '
$ws.Range('D48').Value = "'95.12"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.80%  '
$ws.Range('D49').Value = "'4.72"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.80%  '
$ws.Range('D50').Value = "'66.88"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.25%  '
$ws.Range('D51').Value = "'52.20"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.94%  '
